# Actualización automática 2025-09-26 08:30:09
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M2").Value = 2774.22
$ws1.Range("M11").Value = 30512.6
$ws1.Range("L31").Value = 537.34
$ws1.Range("L60").Value = "3 de 58"
$ws1.Range("M60").Value = "10 de 58"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F2").Value = 2774.22
$ws2.Range("F11").Value = 32451.88
$ws2.Range("F31").Value = 537.34
$ws2.Range("F60").Value = 69656.06

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D9").Value = -10.44
$ws3.Range("E9").Value = 10.44

$ws3.Range("D11").Value = 1607.49
$ws3.Range("E11").Value = 1899.17949822329
$ws3.Range("F11").Value = 0.458409325661988

$ws3.Range("D12").Value = 55050.37
$ws3.Range("E12").Value = -22645.57
$ws3.Range("F12").Value = 1.698833814743495

$ws3.Range("D15").Value = 69822.86000000002
$ws3.Range("E15").Value = -19339.09294897479
$ws3.Range("F15").Value = 1.383075473140273
